$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for the new "Vendors" column
$ws.Range("D1").Value = "Vendors:"

# Vendor assignment per product row (matches the per-vendor breakout sheets)
$ws.Range("D2").Value = "Amazon"
$ws.Range("D3").Value = "Amazon"
$ws.Range("D4").Value = "Amazon"

$ws.Range("D5").Value = "Walmart"
$ws.Range("D6").Value = "Walmart"
$ws.Range("D7").Value = "Walmart"

$ws.Range("D8").Value = "Gamestop"
$ws.Range("D9").Value = "Gamestop"
$ws.Range("D10").Value = "Gamestop"

# Match the saved selection state from the diff
$ws.Range("D9").Select()
